$wb = $excel.ActiveWorkbook

# --- 1. Update the Metadata "Date" value (row 8, column B) ------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2022-03-26T00:24:52-04:00"

# --- 2. Add two new "Include ValueSets" sheets at the end of the workbook ---
$template = $wb.Worksheets.Item("Include ValueSets 3")

# Sheet "Include ValueSets 4" -> LeftRightBothVS
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "Include ValueSets 4"

$template.Range("A1:A2").Copy() | Out-Null
$ws4.Range("A1:A2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws4.Range("A1").Value = "ValueSet URL"
$ws4.Range("A2").Value = "http://hl7.org/fhir/us/pacio-splasch/ValueSet/LeftRightBothVS"

$ws4.Columns.Item(1).ColumnWidth = 30.703125
$ws4.Columns.Item(2).ColumnWidth = 50.703125

# Sheet "Include ValueSets 5" -> HearBetterInOneEarVS
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add($null, $lastSheet2)
$ws5.Name = "Include ValueSets 5"

$template.Range("A1:A2").Copy() | Out-Null
$ws5.Range("A1:A2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws5.Range("A1").Value = "ValueSet URL"
$ws5.Range("A2").Value = "http://hl7.org/fhir/us/pacio-splasch/ValueSet/HearBetterInOneEarVS"

$ws5.Columns.Item(1).ColumnWidth = 30.703125
$ws5.Columns.Item(2).ColumnWidth = 50.703125

# Restore original active sheet/selection
$meta.Activate()
